$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated C-column (temperature) values for the ramping-down schedule
$values = @{
    2  = 9
    3  = 9
    4  = 9
    5  = 9
    6  = 9
    7  = 9
    8  = 12
    9  = 14
    10 = 15
    11 = 15
    12 = 15
    13 = 15
    14 = 15
    15 = 15
    16 = 15
    17 = 15
    18 = 15
    19 = 15
    20 = 15
    21 = 14
    22 = 12
    23 = 9
    24 = 9
    25 = 9
    26 = 9
    27 = 9
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}

# Update the active selection to match the new edit location
$ws.Range("C23").Select()
